# sw/qa/extras/ooxmlexport/data/tscp.docx
#
# DOC(X) filter: BAF -> BAILS in the smart tag mapping.
#
# * Paragraph 1: "before" -> "1st paragraph, non-business." (with "st"
#   as a superscript run).
# * Paragraph 2: drop the stray "_GoBack" bookmark, remap the RDF smart
#   tag's attributes from the old TSCP "BAF" namespace onto the BAILS
#   names, and change "Hello world!" to "2nd paragraph, confidential."
#   (with "nd" superscript).
# * Paragraph 3: "after" -> "3rd paragraph, non-business." (with "rd"
#   superscript).
#
# The smart tag / run-split structure isn't reachable through the
# high-level Word object model in this host, so each paragraph's
# content is replaced wholesale via Range.InsertXML with a literal
# WordprocessingML <w:p> fragment (InsertXML replaces the contents of
# the range it is called on).

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- paragraph 1 ------------------------------------------------------
$p1Xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:rPr/>
    <w:t>1</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>st</w:t>
  </w:r>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"> paragraph, non-business.</w:t>
  </w:r>
</w:p>
"@
$d.Paragraphs.Item(1).Range.InsertXML($p1Xml)

# --- paragraph 2 (smart tag BAF -> BAILS) -----------------------------
$p2Xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:rPr/>
  </w:pPr>
  <w:smartTag w:uri="http://www.w3.org/1999/02/22-rdf-syntax-ns#" w:element="RDF">
    <w:smartTagPr>
      <w:attr w:name="urn:bails:ExportControl:Authorization:StartValidity" w:val="2015-11-27"/>
      <w:attr w:name="urn:bails:ExportControl:BusinessAuthorization:Identifier" w:val="urn:example:tscp:1"/>
      <w:attr w:name="urn:bails:ExportControl:BusinessAuthorizationCategory:Identifier" w:val="urn:example:tscp:1:confidential"/>
    </w:smartTagPr>
  </w:smartTag>
  <w:r>
    <w:rPr/>
    <w:t>2</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>nd</w:t>
  </w:r>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"> paragraph, confidential.</w:t>
  </w:r>
</w:p>
"@
$d.Paragraphs.Item(2).Range.InsertXML($p2Xml)

# --- paragraph 3 -------------------------------------------------------
$p3Xml = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:rPr/>
  </w:pPr>
  <w:r>
    <w:rPr/>
    <w:t>3</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:vertAlign w:val="superscript"/>
    </w:rPr>
    <w:t>rd</w:t>
  </w:r>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"> paragraph, non-business.</w:t>
  </w:r>
</w:p>
"@
$d.Paragraphs.Item(3).Range.InsertXML($p3Xml)

Write-Output "paragraph 1: $($d.Paragraphs.Item(1).Range.Text)"
Write-Output "paragraph 2: $($d.Paragraphs.Item(2).Range.Text)"
Write-Output "paragraph 3: $($d.Paragraphs.Item(3).Range.Text)"
